# Update the date heading and the 25 multiplication problems in the table.
$d = $word.ActiveDocument

# 1) Date heading: 2024-09-15 Sunday -> 2024-09-16 Monday
$d.Content.Find.Execute("2024-09-15 Sunday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-09-16 Monday", 2)

# 2) Multiplication table (5 data rows at table rows 1, 5, 10, 15, 20 x 5 columns).
$t = $d.Tables.Item(1)

$newValues = @{
    1  = @("18×22=", "95×53=", "55×71=", "43×51=", "56×74=")
    5  = @("62×66=", "36×25=", "54×25=", "90×52=", "95×53=")
    10 = @("82×92=", "90×99=", "27×78=", "53×86=", "12×86=")
    15 = @("40×66=", "22×30=", "99×86=", "25×33=", "66×33=")
    20 = @("96×18=", "54×92=", "85×95=", "53×66=", "21×75=")
}

foreach ($rowIndex in $newValues.Keys) {
    $values = $newValues[$rowIndex]
    for ($col = 1; $col -le $values.Count; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
